$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells to Text format before writing so that values
# like "285.98" are stored as literal strings (matching the original
# inlineStr/shared-string cells) instead of being auto-converted to numbers.
$dRange = $ws.Range("D2:D50")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "22.396.72"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "1.566.02"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").Value = "285.98"
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").Value = "0.3703"
$ws.Range("E7").Value = "  +1.65%  "
$ws.Range("D8").Value = "0.3279"
$ws.Range("E8").Value = "  -1.55%  "
$ws.Range("D9").Value = "46.43"
$ws.Range("E9").Value = "  -4.57%  "
$ws.Range("D10").Value = "1.148"
$ws.Range("E10").Value = "  +2.12%  "
$ws.Range("D11").Value = "0.07409"
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").Value = "20.45"
$ws.Range("E13").Value = "  -1.58%  "
$ws.Range("D14").Value = "5.839"
$ws.Range("E14").Value = "  -1.99%  "
$ws.Range("D15").Value = "6.812"
$ws.Range("E15").Value = "  -1.36%  "
$ws.Range("D16").Value = "1.562.00"
$ws.Range("E16").Value = "  -0.74%  "
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("D18").Value = "0.06692"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("D19").Value = "86.09"
$ws.Range("E19").Value = "  -2.36%  "
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("D21").Value = "6.333"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").Value = "16.28"
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("D23").Value = "11.80"
$ws.Range("E23").Value = "  -1.80%  "
$ws.Range("D24").Value = "22.401.29"
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("D25").Value = "2.321"
$ws.Range("E25").Value = "  -2.59%  "
$ws.Range("D26").Value = "2.584"
$ws.Range("E26").Value = "  +1.77%  "
$ws.Range("D27").Value = "150.65"
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("D28").Value = "19.31"
$ws.Range("E28").Value = "  -0.42%  "
$ws.Range("D29").Value = "4.944"
$ws.Range("E29").Value = "  -1.08%  "
$ws.Range("D30").Value = "123.49"
$ws.Range("E30").Value = "  -0.40%  "
$ws.Range("D31").Value = "1.740.80"
$ws.Range("E31").Value = "  -1.09%  "
$ws.Range("D32").Value = "1.048"
$ws.Range("E32").Value = "  +0.33%  "
$ws.Range("D33").Value = "1.971"
$ws.Range("E33").Value = "  -0.99%  "
$ws.Range("D34").Value = "5.978"
$ws.Range("E34").Value = "  -2.12%  "
$ws.Range("D35").Value = "9.697"
$ws.Range("E35").Value = "  -1.09%  "
$ws.Range("D36").Value = "0.08247"
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("D37").Value = "0.02395"
$ws.Range("E37").Value = "  -0.86%  "
$ws.Range("D38").Value = "1.300"
$ws.Range("E38").Value = "  +0.92%  "
$ws.Range("D39").Value = "0.06331"
$ws.Range("E39").Value = "  -1.63%  "
$ws.Range("D40").Value = "0.2189"
$ws.Range("E40").Value = "  -2.07%  "
$ws.Range("D41").Value = "5.228"
$ws.Range("E41").Value = "  -2.46%  "
$ws.Range("D42").Value = "11.15"
$ws.Range("E42").Value = "  -0.50%  "
$ws.Range("D43").Value = "0.6121"
$ws.Range("E43").Value = "  -2.10%  "
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("D45").Value = "13.76"
$ws.Range("E45").Value = "  -0.71%  "
$ws.Range("D46").Value = "0.5955"
$ws.Range("E46").Value = "  -1.56%  "
$ws.Range("D47").Value = "3.746"
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("D48").Value = "2.014"
$ws.Range("E48").Value = "  -1.13%  "
$ws.Range("D49").Value = "123.81"
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("D50").Value = "1.181"
$ws.Range("E50").Value = "  -2.67%  "

# Restore the default (unstyled) cell format on the Price column so the
# resulting XML matches the original cells that carry no explicit style.
$dRange.Style = "Normal"

